$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ausruestung / Attribute aktualisieren (Inventar einlesen, Waffe/Ruestung ausruesten, Trefferprobe)
$ws.Range("B7").Value = 14   # MU
$ws.Range("B8").Value = 14   # KL
$ws.Range("B9").Value = 10   # IN
$ws.Range("B10").Value = 10  # CH
$ws.Range("B13").Value = 10  # KO

# "LeP" -> "LP"
$ws.Range("A16").Value = "LP"

# Auswahl wie im Original-Commit auf A16
[void]$ws.Range("A16").Select()
